# Fix status names (shared-string text corrections):
#   "bleu" -> "noir"
#   "résultat et / ou publication posté" -> "résultat postés ou publiés"
#   "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
#   "résultat et / ou publication posté dans les 36 mois" -> "résultat postés ou publiés dans les 36 mois"
#   "résultat et / ou publication posté dans les 12 mois" -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "bleu" = "noir"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
}

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $map.ContainsKey($val)) {
            $cell.Value2 = $map[$val]
        }
    }
}
